$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.206.38"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.163.22"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "253.11"
$ws.Range("E5").Value = "  +6.54%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").Value = "72.95"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "39.62"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "6.72"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "2.486.40"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "14.15"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "2.142.34"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "0.763"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "42.052.03"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "70.38"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "5.81"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "225.92"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "9.53"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("E24").Value = "  +6.26%  "
$ws.Range("D26").Value = "10.42"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").Value = "36.40"
$ws.Range("E30").Value = "  +10.71%  "
$ws.Range("D31").Value = "168.14"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "19.91"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "0.0803"
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  +5.15%  "
$ws.Range("D37").Value = "4.23"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "0.0330"
$ws.Range("E38").Value = "  +8.93%  "
$ws.Range("D39").Value = "11.79"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "0.195"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D42").Value = "58.51"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "5.11"
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").Value = "102.16"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("B45").Value = "WOONetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D45").Value = "0.460"
$ws.Range("E45").Value = "  +15.67%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.19"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("D47").Value = "0.0962"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +9.46%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "2.64"
$ws.Range("E51").Value = "  +0.79%  "
